$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous "NA" value in C122 (becomes an empty text cell),
# matching the style of an existing empty text cell so no stray format
# is introduced.
$ws.Range("C122").Formula = "'"
$ws.Range("C122").Style = $ws.Range("C2").Style

# Append the new row 123 with the data that used to live in row 122's
# "NA" slot, now dated 2025-06-03.
$ws.Range("A123").Formula = "'2025-06-03"
$ws.Range("A123").Style = $ws.Range("A122").Style

$ws.Range("B123").Value2 = "Rien ne nous concerne aujourd'hui !"
$ws.Range("B123").Style = $ws.Range("B122").Style

$ws.Range("C123").Formula = "'NA"
$ws.Range("C123").Style = $ws.Range("C121").Style

$ws.Range("D123").Value2 = 1
$ws.Range("D123").Style = $ws.Range("D122").Style
